$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -13.376
$ws.Range("B9").Value = 5.496
$ws.Range("C9").Value = -11.122
$ws.Range("D9").Value = -7.244
$ws.Range("C11").Value = -11.802
$ws.Range("B13").Value = 6.017999999999999
$ws.Range("B16").Value = 5.977
$ws.Range("C16").Value = -13.234
$ws.Range("B18").Value = 5.257
$ws.Range("B20").Value = 6.827
$ws.Range("D22").Value = -8.121
$ws.Range("C23").Value = -13.331
$ws.Range("C24").Value = -12.424
$ws.Range("B26").Value = 5.077
$ws.Range("C26").Value = -11.559
$ws.Range("B27").Value = 6.412999999999999
$ws.Range("D27").Value = -7.866
$ws.Range("B29").Value = 5.265
$ws.Range("D29").Value = -7.568000000000001
$ws.Range("D32").Value = -7.252
$ws.Range("C34").Value = -12.105
$ws.Range("B35").Value = 8.597999999999999
$ws.Range("C35").Value = -12.578
$ws.Range("B36").Value = 8.043999999999999
$ws.Range("D37").Value = -7.715000000000001
$ws.Range("D38").Value = -7.833
$ws.Range("D39").Value = -7.56
$ws.Range("D41").Value = -8.15
$ws.Range("C44").Value = -12.798
$ws.Range("B45").Value = 5.994000000000001
$ws.Range("D45").Value = -8.153
$ws.Range("C48").Value = -12.478
$ws.Range("D48").Value = -7.539
$ws.Range("C49").Value = -12.687
$ws.Range("D51").Value = -8.174000000000001
$ws.Range("C52").Value = -11.717
$ws.Range("B55").Value = 5.278
$ws.Range("D56").Value = -8.266
$ws.Range("B57").Value = 5.116
$ws.Range("D57").Value = -8.31
$ws.Range("D61").Value = -7.900999999999999
$ws.Range("D64").Value = -7.785999999999999
$ws.Range("C66").Value = -11.278
$ws.Range("C67").Value = -11.547
$ws.Range("B69").Value = 5.362
$ws.Range("C73").Value = -11.5
$ws.Range("D75").Value = -8.068999999999999
$ws.Range("B76").Value = 6.182
$ws.Range("B78").Value = 8.891
$ws.Range("C78").Value = -11.793
$ws.Range("C80").Value = -11.83
$ws.Range("B82").Value = 5.961
$ws.Range("D82").Value = -8.545999999999999
$ws.Range("B83").Value = 5.906000000000001
$ws.Range("D90").Value = -7.306
$ws.Range("C91").Value = -13.133
$ws.Range("B93").Value = 5.008
$ws.Range("D93").Value = -7.008
$ws.Range("B97").Value = 5.476
$ws.Range("C97").Value = -10.517
$ws.Range("C99").Value = -11.4
$ws.Range("D102").Value = -7.471000000000001
$ws.Range("C104").Value = -13.31
$ws.Range("D105").Value = -7.536000000000001
